$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) from row 4 down into row 5 for the two styled columns
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("G4").Copy()
$ws.Range("G5").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A5").Value = 42647.681655092594
$ws.Range("B5").Value = $false
$ws.Range("C5").Value = 9931.74
$ws.Range("D5").Value = 10029.02
$ws.Range("E5").Value = 313
$ws.Range("F5").Value = 309.97000000000003
$ws.Range("G5").Value = $false
$ws.Range("H5").Value = -0.97
$ws.Range("I5").Value = $true
